$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.931133031845093
$ws.Range("B1").Value = 5.575462341308594
$ws.Range("C1").Value = 4.579350471496582
$ws.Range("D1").Value = 5.338911533355713
$ws.Range("E1").Value = 4.508998870849609
